{"js": "// Applies the README.md report edits described by the commit:\n//   - expands the \"progession\" bullet with the new-species / 12-fish sentence\n//   - adds \"c\u1ed1t truy\u1ec7n\" to the \"M\u1edf r\u1ed9ng n\u1ed9i dung\" bullet\n//   - expands \"game c\u00f3 t\u00ean tu\u1ed5i\" -> \"game E-sport c\u00f3 t\u00ean tu\u1ed5i \u0111ang th\u1ed1ng tr\u1ecb l\u00e0ng game\"\n//   - adds \"d\u1ec5 t\u00ednh v\u00e0 th\u00edch s\u1ef1 \u0111\u01a1n gi\u1ea3n\" reasoning to \"H\u01b0\u1edbng gi\u1ea3i quy\u1ebft\"\n//   - adds \"theo quy t\u1eafc chung\" to the variable naming convention sentence\n//   - replaces \"b\u1ecf d\u1edf\" with \"b\u1ee1 ng\u1ee1\" in the gameplay-difficulty bullet\n//\n// Each replacement is done via body.search(...) on the exact original text\n// (so it is a no-op / safe if the text was already changed) followed by\n// insertText(..., \"Replace\") on the located range.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    return false;\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n  return true;\n}\n\n// 1) \"C\u1ea3m gi\u00e1c ph\u00e1t tri\u1ec3n theo th\u1eddi gian (progession)\" bullet: insert the new\n//    sentence about later levels / new fish species / 12 species total\n//    before the closing \"\u0110\u00e2y ch\u00ednh l\u00e0 y\u1ebfu t\u1ed1...\" sentence.\nawait replaceOnce(\n  \" \u0110\u00e2y ch\u00ednh l\u00e0 y\u1ebfu t\u1ed1 m\u00e0 nhi\u1ec1u game casual tr\u00ean Google play \u0111ang thi\u1ebfu\",\n  \" Ngo\u00e0i ra, nh\u1eefng level sau Player s\u1ebd \u0111\u01b0\u1ee3c ti\u1ebfp c\u1eadn v\u1edbi nh\u1eefng lo\u00e0i c\u00e1 m\u1edbi\" +\n    \"(trong game c\u00f3 t\u1ea5t c\u1ea3 12 lo\u00e0i c\u00e1). \" +\n    \"\u0110\u00e2y ch\u00ednh l\u00e0 y\u1ebfu t\u1ed1 m\u00e0 nhi\u1ec1u game casual tr\u00ean Google play \u0111ang thi\u1ebfu\"\n);\n\n// 2) \"M\u1edf r\u1ed9ng n\u1ed9i dung\" bullet: add \"c\u1ed1t truy\u1ec7n\" to the list of things that\n//    could be added to the game.\nawait replaceOnce(\n  \" Hi\u1ec7n t\u1ea1i game c\u00f3 s\u1ed1 l\u01b0\u1ee3ng lo\u00e0i c\u00e1 v\u00e0 m\u00f4i tr\u01b0\u1eddng h\u1ea1n ch\u1ebf. Vi\u1ec7c b\u1ed5 sung th\u00eam nhi\u1ec1u lo\u1ea1i sinh v\u1eadt bi\u1ec3n, nhi\u1ec7m v\u1ee5, v\u00e0 ch\u1ebf \u0111\u1ed9 ch\u01a1i s\u1ebd gi\u00fap t\u0103ng s\u1ef1 h\u1ea5p d\u1eabn\",\n  \" Hi\u1ec7n t\u1ea1i game c\u00f3 s\u1ed1 l\u01b0\u1ee3ng lo\u00e0i c\u00e1 v\u00e0 m\u00f4i tr\u01b0\u1eddng h\u1ea1n ch\u1ebf. Vi\u1ec7c b\u1ed5 sung th\u00eam nhi\u1ec1u lo\u1ea1i sinh v\u1eadt bi\u1ec3n, nhi\u1ec7m v\u1ee5, c\u1ed1t truy\u1ec7n v\u00e0 ch\u1ebf \u0111\u1ed9 ch\u01a1i s\u1ebd gi\u00fap t\u0103ng s\u1ef1 h\u1ea5p d\u1eabn\"\n);\n\n// 3) \"Thi\u1ebfu kinh ph\u00ed \u0111\u1ec3 ph\u00e1t tri\u1ec3n\" bullet: \"game c\u00f3 t\u00ean tu\u1ed5i\" ->\n//    \"game E-sport c\u00f3 t\u00ean tu\u1ed5i \u0111ang th\u1ed1ng tr\u1ecb l\u00e0ng game\" (first occurrence\n//    only \u2014 the phrase also appears later in \"H\u01b0\u1edbng gi\u1ea3i quy\u1ebft\").\nawait replaceOnce(\n  \"game c\u00f3 t\u00ean tu\u1ed5i\",\n  \"game E-sport c\u00f3 t\u00ean tu\u1ed5i \u0111ang th\u1ed1ng tr\u1ecb l\u00e0ng game\"\n);\n\n// 4) \"H\u01b0\u1edbng gi\u1ea3i quy\u1ebft\" bullet: explain *why* young players are targeted.\nawait replaceOnce(\n  \" : Game s\u1ebd t\u1eadp trung v\u00e0o nh\u1eefng \u0111\u1ed1i t\u01b0\u1ee3ng nh\u1ecf tu\u1ed5i, b\u1edfi v\u00ec \u0111\u00e2y l\u00e0 \u0111\u1ed1i t\u01b0\u1ee3ng c\u00f2n h\u1ee9ng th\u00fa v\u1edbi th\u1ebf gi\u1edbi b\u00ean ngo\u00e0i. Nh\u1eefng b\u1ea1n l\u1edbn tu\u1ed5i \u0111a ph\u1ea7n s\u1ebd trung th\u00e0nh v\u1edbi c\u00e1c game c\u00f3 t\u00ean tu\u1ed5i ho\u1eb7c l\u00e0 triple A\",\n  \" : Game s\u1ebd t\u1eadp trung v\u00e0o nh\u1eefng \u0111\u1ed1i t\u01b0\u1ee3ng nh\u1ecf tu\u1ed5i, b\u1edfi v\u00ec \u0111\u00e2y l\u00e0 \u0111\u1ed1i t\u01b0\u1ee3ng d\u1ec5 t\u00ednh v\u00e0 th\u00edch s\u1ef1 \u0111\u01a1n gi\u1ea3n. Nh\u1eefng b\u1ea1n l\u1edbn tu\u1ed5i \u0111a ph\u1ea7n s\u1ebd trung th\u00e0nh v\u1edbi c\u00e1c game c\u00f3 t\u00ean tu\u1ed5i ho\u1eb7c l\u00e0 triple A\"\n);\n\n// 5) \"Clean code\" bullet: variable naming convention sentence.\nawait replaceOnce(\n  \"Quy \u01b0\u1edbc \u0111\u1eb7t t\u00ean bi\u1ebfn \u0111\u1ec3 t\u1ea5t c\u1ea3 m\u1ecdi ng\u01b0\u1eddi c\u00f9ng hi\u1ec3u. \",\n  \"Quy \u01b0\u1edbc \u0111\u1eb7t t\u00ean bi\u1ebfn theo quy t\u1eafc chung \u0111\u1ec3 t\u1ea5t c\u1ea3 m\u1ecdi ng\u01b0\u1eddi c\u00f9ng hi\u1ec3u. \"\n);\n\n// 6) \"Thi\u1ebft k\u1ebf gameplay\" bullet: \"b\u1ecf d\u1edf\" -> \"b\u1ee1 ng\u1ee1\".\nawait replaceOnce(\n  \", v\u00e0 c\u0169ng kh\u00f4ng t\u0103ng \u0111\u1ed9t ng\u1ed9t khi\u1ebfn ng\u01b0\u1eddi ch\u01a1i ph\u1ea3i b\u1ecf d\u1edf\",\n  \", v\u00e0 c\u0169ng kh\u00f4ng t\u0103ng \u0111\u1ed9t ng\u1ed9t khi\u1ebfn ng\u01b0\u1eddi ch\u01a1i ph\u1ea3i b\u1ee1 ng\u1ee1\"\n);\n", "ps1": "# Applies the README.md report edits described by the commit:\n#   - expands the \"progession\" bullet with the new-species / 12-fish sentence\n#   - adds \"c\u1ed1t truy\u1ec7n\" to the \"M\u1edf r\u1ed9ng n\u1ed9i dung\" bullet\n#   - expands \"game c\u00f3 t\u00ean tu\u1ed5i\" -> \"game E-sport c\u00f3 t\u00ean tu\u1ed5i \u0111ang th\u1ed1ng tr\u1ecb l\u00e0ng game\"\n#   - adds \"d\u1ec5 t\u00ednh v\u00e0 th\u00edch s\u1ef1 \u0111\u01a1n gi\u1ea3n\" reasoning to \"H\u01b0\u1edbng gi\u1ea3i quy\u1ebft\"\n#   - adds \"theo quy t\u1eafc chung\" to the variable naming convention sentence\n#   - replaces \"b\u1ecf d\u1edf\" with \"b\u1ee1 ng\u1ee1\" in the gameplay-difficulty bullet\n#\n# Each replacement uses Find/Replace (wdReplaceOne) scoped to the exact\n# original sentence, so it only ever touches the single intended run/\n# occurrence (duplicated substrings elsewhere in the doc, e.g. the second\n# \"game c\u00f3 t\u00ean tu\u1ed5i\", are left untouched).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\n# 1) \"C\u1ea3m gi\u00e1c ph\u00e1t tri\u1ec3n theo th\u1eddi gian (progession)\" bullet: insert the new\n#    sentence about later levels / new fish species / 12 species total\n#    before the closing \"\u0110\u00e2y ch\u00ednh l\u00e0 y\u1ebfu t\u1ed1...\" sentence.\nReplace-Once \" \u0110\u00e2y ch\u00ednh l\u00e0 y\u1ebfu t\u1ed1 m\u00e0 nhi\u1ec1u game casual tr\u00ean Google play \u0111ang thi\u1ebfu\" \" Ngo\u00e0i ra, nh\u1eefng level sau Player s\u1ebd \u0111\u01b0\u1ee3c ti\u1ebfp c\u1eadn v\u1edbi nh\u1eefng lo\u00e0i c\u00e1 m\u1edbi(trong game c\u00f3 t\u1ea5t c\u1ea3 12 lo\u00e0i c\u00e1). \u0110\u00e2y ch\u00ednh l\u00e0 y\u1ebfu t\u1ed1 m\u00e0 nhi\u1ec1u game casual tr\u00ean Google play \u0111ang thi\u1ebfu\"\n\n# 2) \"M\u1edf r\u1ed9ng n\u1ed9i dung\" bullet: add \"c\u1ed1t truy\u1ec7n\" to the list of things that\n#    could be added to the game.\nReplace-Once \" Hi\u1ec7n t\u1ea1i game c\u00f3 s\u1ed1 l\u01b0\u1ee3ng lo\u00e0i c\u00e1 v\u00e0 m\u00f4i tr\u01b0\u1eddng h\u1ea1n ch\u1ebf. Vi\u1ec7c b\u1ed5 sung th\u00eam nhi\u1ec1u lo\u1ea1i sinh v\u1eadt bi\u1ec3n, nhi\u1ec7m v\u1ee5, v\u00e0 ch\u1ebf \u0111\u1ed9 ch\u01a1i s\u1ebd gi\u00fap t\u0103ng s\u1ef1 h\u1ea5p d\u1eabn\" \" Hi\u1ec7n t\u1ea1i game c\u00f3 s\u1ed1 l\u01b0\u1ee3ng lo\u00e0i c\u00e1 v\u00e0 m\u00f4i tr\u01b0\u1eddng h\u1ea1n ch\u1ebf. Vi\u1ec7c b\u1ed5 sung th\u00eam nhi\u1ec1u lo\u1ea1i sinh v\u1eadt bi\u1ec3n, nhi\u1ec7m v\u1ee5, c\u1ed1t truy\u1ec7n v\u00e0 ch\u1ebf \u0111\u1ed9 ch\u01a1i s\u1ebd gi\u00fap t\u0103ng s\u1ef1 h\u1ea5p d\u1eabn\"\n\n# 3) \"Thi\u1ebfu kinh ph\u00ed \u0111\u1ec3 ph\u00e1t tri\u1ec3n\" bullet: \"game c\u00f3 t\u00ean tu\u1ed5i\" ->\n#    \"game E-sport c\u00f3 t\u00ean tu\u1ed5i \u0111ang th\u1ed1ng tr\u1ecb l\u00e0ng game\" (wdReplaceOne hits\n#    only the first occurrence \u2014 the phrase also appears later, unchanged,\n#    in \"H\u01b0\u1edbng gi\u1ea3i quy\u1ebft\").\nReplace-Once \"game c\u00f3 t\u00ean tu\u1ed5i\" \"game E-sport c\u00f3 t\u00ean tu\u1ed5i \u0111ang th\u1ed1ng tr\u1ecb l\u00e0ng game\"\n\n# 4) \"H\u01b0\u1edbng gi\u1ea3i quy\u1ebft\" bullet: explain *why* young players are targeted.\nReplace-Once \" : Game s\u1ebd t\u1eadp trung v\u00e0o nh\u1eefng \u0111\u1ed1i t\u01b0\u1ee3ng nh\u1ecf tu\u1ed5i, b\u1edfi v\u00ec \u0111\u00e2y l\u00e0 \u0111\u1ed1i t\u01b0\u1ee3ng c\u00f2n h\u1ee9ng th\u00fa v\u1edbi th\u1ebf gi\u1edbi b\u00ean ngo\u00e0i. Nh\u1eefng b\u1ea1n l\u1edbn tu\u1ed5i \u0111a ph\u1ea7n s\u1ebd trung th\u00e0nh v\u1edbi c\u00e1c game c\u00f3 t\u00ean tu\u1ed5i ho\u1eb7c l\u00e0 triple A\" \" : Game s\u1ebd t\u1eadp trung v\u00e0o nh\u1eefng \u0111\u1ed1i t\u01b0\u1ee3ng nh\u1ecf tu\u1ed5i, b\u1edfi v\u00ec \u0111\u00e2y l\u00e0 \u0111\u1ed1i t\u01b0\u1ee3ng d\u1ec5 t\u00ednh v\u00e0 th\u00edch s\u1ef1 \u0111\u01a1n gi\u1ea3n. Nh\u1eefng b\u1ea1n l\u1edbn tu\u1ed5i \u0111a ph\u1ea7n s\u1ebd trung th\u00e0nh v\u1edbi c\u00e1c game c\u00f3 t\u00ean tu\u1ed5i ho\u1eb7c l\u00e0 triple A\"\n\n# 5) \"Clean code\" bullet: variable naming convention sentence.\nReplace-Once \"Quy \u01b0\u1edbc \u0111\u1eb7t t\u00ean bi\u1ebfn \u0111\u1ec3 t\u1ea5t c\u1ea3 m\u1ecdi ng\u01b0\u1eddi c\u00f9ng hi\u1ec3u. \" \"Quy \u01b0\u1edbc \u0111\u1eb7t t\u00ean bi\u1ebfn theo quy t\u1eafc chung \u0111\u1ec3 t\u1ea5t c\u1ea3 m\u1ecdi ng\u01b0\u1eddi c\u00f9ng hi\u1ec3u. \"\n\n# 6) \"Thi\u1ebft k\u1ebf gameplay\" bullet: \"b\u1ecf d\u1edf\" -> \"b\u1ee1 ng\u1ee1\".\nReplace-Once \", v\u00e0 c\u0169ng kh\u00f4ng t\u0103ng \u0111\u1ed9t ng\u1ed9t khi\u1ebfn ng\u01b0\u1eddi ch\u01a1i ph\u1ea3i b\u1ecf d\u1edf\" \", v\u00e0 c\u0169ng kh\u00f4ng t\u0103ng \u0111\u1ed9t ng\u1ed9t khi\u1ebfn ng\u01b0\u1eddi ch\u01a1i ph\u1ea3i b\u1ee1 ng\u1ee1\"\n"}
